# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) listed periods 2207..2112 in
# descending (most-recent-first) order. The update re-sorts them into
# chronological (ascending) order 2112..2207, and refreshes the
# "Valor Mora" (F) / "Salario Basico" (G) figures that go with the
# re-sorted periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological period order for rows 16..23.
$periods = @("2112", "2201", "2202", "2203", "2204", "2205", "2206", "2207")

# Valor Mora (F) for each of those periods - every period is a full
# month (40000) except the most recent one (2207), which is a partial
# period (33333), matching the same pattern the sheet had before
# (the most-recent partial period carried the 33333 value).
$valorMora = @{
    "2112" = 40000
    "2201" = 40000
    "2202" = 40000
    "2203" = 40000
    "2204" = 40000
    "2205" = 40000
    "2206" = 40000
    "2207" = 33333
}

# Salario Basico (G) refreshed uniformly for every period.
$salarioBasico = 900000

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value2 = $period
    $ws.Cells.Item($row, 6).Value2 = $valorMora[$period]
    $ws.Cells.Item($row, 7).Value2 = $salarioBasico
}
